# New weekly price record for "Bruselas (repollito)" needs to be inserted
# right after the existing row 61, pushing all subsequent rows down by one
# (old row 62 becomes row 63, ..., old row 103 becomes row 104).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 62 - shifts rows 62:103 down to 63:104.
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new record. Columns that are
# constant across every data row in this sheet (A, B, C, E, F, G, H, I, N,
# O, Q, R) are copied from the surrounding rows; the remaining columns hold
# the new observation's values.
$ws.Range("A62").Value = 10
$ws.Range("B62").Value = "Vega Modelo de Temuco"
$ws.Range("C62").Value = "La Araucanía"
$ws.Range("D62").Value = 44767
$ws.Range("E62").Value = 9
$ws.Range("F62").Value = 100112035
$ws.Range("G62").Value = "Bruselas (repollito)"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 200
$ws.Range("K62").Value = 25000
$ws.Range("L62").Value = 26000
$ws.Range("M62").Value = 25500
$ws.Range("N62").Value = "$/malla 10 kilos"
$ws.Range("O62").Value = "Provincia de Quillota"
$ws.Range("P62").Value = 2550
$ws.Range("Q62").Value = 10
$ws.Range("R62").Value = "Hortaliza"
